$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2024-12-28 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-12-29 Sunday", 2) | Out-Null

# Update the practice-problem table. Data rows are Word table rows 1,5,9,13,17
# (rows 2-4, 6-8, 10-12, 14-16, 18-20 are blank spacer rows).
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "83÷6=13, 5"
$t.Cell(1, 2).Range.Text = "49÷8=6, 1"
$t.Cell(1, 3).Range.Text = "43÷5=8, 3"
$t.Cell(1, 4).Range.Text = "93÷6=15, 3"
$t.Cell(1, 5).Range.Text = "50÷6=8, 2"
$t.Cell(5, 1).Range.Text = "34÷9=3, 7"
$t.Cell(5, 2).Range.Text = "88÷4=22, 0"
$t.Cell(5, 3).Range.Text = "86÷7=12, 2"
$t.Cell(5, 4).Range.Text = "76÷9=8, 4"
$t.Cell(5, 5).Range.Text = "29÷4=7, 1"
$t.Cell(9, 1).Range.Text = "90÷9=10, 0"
$t.Cell(9, 2).Range.Text = "97÷2=48, 1"
$t.Cell(9, 3).Range.Text = "84÷7=12, 0"
$t.Cell(9, 4).Range.Text = "53÷3=17, 2"
$t.Cell(9, 5).Range.Text = "78÷6=13, 0"
$t.Cell(13, 1).Range.Text = "80÷9=8, 8"
$t.Cell(13, 2).Range.Text = "37÷6=6, 1"
$t.Cell(13, 3).Range.Text = "13÷5=2, 3"
$t.Cell(13, 4).Range.Text = "82÷5=16, 2"
$t.Cell(13, 5).Range.Text = "82÷6=13, 4"
$t.Cell(17, 1).Range.Text = "74÷6=12, 2"
$t.Cell(17, 2).Range.Text = "82÷9=9, 1"
$t.Cell(17, 3).Range.Text = "81÷7=11, 4"
$t.Cell(17, 4).Range.Text = "91÷9=10, 1"
$t.Cell(17, 5).Range.Text = "71÷5=14, 1"
